# Daily attendance processing - 2025-11-08 05:45:23
#
# Normalize the "Recorded By" (column G) entries on the session analysis
# sheet: wherever the literal "System" recorder appears first in the
# comma-separated list alongside another recorder, move it so it is no
# longer listed first (fixes a sort/display quirk from the importer).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Recorded By") runs from row 2 down to the last populated row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row   # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Text

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, admin@admin.com") {
        $cell.Value = "admin@admin.com, System"
    }
    elseif ($val -eq "System, system, backup@backdoor.com") {
        $cell.Value = "system, System, backup@backdoor.com"
    }
}
